# Update the pl_mw.xlsx results sheet ("case with 380 kV done") with the
# newly computed line power-flow values for rows 2-25 (columns B,C,E,F,G,H,I,K,L,M,O).
# Columns A, D, J, N are untouched (index / zero columns), as is the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.390635051968502
$ws.Cells.Item(2, 3).Value = 0.2097675727620683
$ws.Cells.Item(2, 5).Value = 0.1381766316326338
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.5515911136399261
$ws.Cells.Item(2, 8).Value = 0.7069331185372363
$ws.Cells.Item(2, 9).Value = 0.7335909157176328
$ws.Cells.Item(2, 11).Value = 0.2662655900770687
$ws.Cells.Item(2, 12).Value = 0.2052607142353935
$ws.Cells.Item(2, 13).Value = 0.1256452696574684
$ws.Cells.Item(2, 15).Value = 2.490416174936882
$ws.Cells.Item(3, 2).Value = 0.3519467378087882
$ws.Cells.Item(3, 3).Value = 0.2098208785555045
$ws.Cells.Item(3, 5).Value = 0.1389951105589731
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.5590628944154794
$ws.Cells.Item(3, 8).Value = 0.7138658593186058
$ws.Cells.Item(3, 9).Value = 0.7428042633263772
$ws.Cells.Item(3, 11).Value = 0.2324442168110181
$ws.Cells.Item(3, 12).Value = 0.202698743479246
$ws.Cells.Item(3, 13).Value = 0.1182093900631536
$ws.Cells.Item(3, 15).Value = 2.520561902114025
$ws.Cells.Item(4, 2).Value = 0.3281787753929564
$ws.Cells.Item(4, 3).Value = 0.2098878727125708
$ws.Cells.Item(4, 5).Value = 0.1395581039834148
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.5640075442404715
$ws.Cells.Item(4, 8).Value = 0.7184011555924705
$ws.Cells.Item(4, 9).Value = 0.7488211929117128
$ws.Cells.Item(4, 11).Value = 0.211589511569187
$ws.Cells.Item(4, 12).Value = 0.2012236513015893
$ws.Cells.Item(4, 13).Value = 0.1136731010964169
$ws.Cells.Item(4, 15).Value = 2.540404361723894
$ws.Cells.Item(5, 2).Value = 0.3184905846082415
$ws.Cells.Item(5, 3).Value = 0.2099238335314979
$ws.Cells.Item(5, 5).Value = 0.1398027523298904
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.566112248966796
$ws.Cells.Item(5, 8).Value = 0.7203194254380314
$ws.Cells.Item(5, 9).Value = 0.7513636592680868
$ws.Cells.Item(5, 11).Value = 0.2030694780222859
$ws.Cells.Item(5, 12).Value = 0.2006472604317864
$ws.Cells.Item(5, 13).Value = 0.1118320662169445
$ws.Cells.Item(5, 15).Value = 2.54882550974488
$ws.Cells.Item(6, 2).Value = 0.3168817361411698
$ws.Cells.Item(6, 3).Value = 0.2099303291511276
$ws.Cells.Item(6, 5).Value = 0.13984429621887
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.566467151339257
$ws.Cells.Item(6, 8).Value = 0.7206421880768588
$ws.Cells.Item(6, 9).Value = 0.7517913015550288
$ws.Cells.Item(6, 11).Value = 0.2016534482787335
$ws.Cells.Item(6, 12).Value = 0.2005530468005361
$ws.Cells.Item(6, 13).Value = 0.1115268236941986
$ws.Cells.Item(6, 15).Value = 2.550244080041786
$ws.Cells.Item(7, 2).Value = 0.3280481263908825
$ws.Cells.Item(7, 3).Value = 0.2098883225695332
$ws.Cells.Item(7, 5).Value = 0.1395613417222474
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.5640355658016212
$ws.Cells.Item(7, 8).Value = 0.7184267421491768
$ws.Cells.Item(7, 9).Value = 0.7488551149685421
$ws.Cells.Item(7, 11).Value = 0.2114746939405592
$ws.Cells.Item(7, 12).Value = 0.2012157776827266
$ws.Cells.Item(7, 13).Value = 0.1136482415237268
$ws.Cells.Item(7, 15).Value = 2.540516575085285
$ws.Cells.Item(8, 2).Value = 0.3772985013144989
$ws.Cells.Item(8, 3).Value = 0.209778867294709
$ws.Cells.Item(8, 5).Value = 0.1384463128200277
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.5540932776346565
$ws.Cells.Item(8, 8).Value = 0.7092657432005218
$ws.Cells.Item(8, 9).Value = 0.7366929940569875
$ws.Cells.Item(8, 11).Value = 0.2546226562072036
$ws.Cells.Item(8, 12).Value = 0.2043570639109049
$ws.Cells.Item(8, 13).Value = 0.1230753723292146
$ws.Cells.Item(8, 15).Value = 2.500533804035285
$ws.Cells.Item(9, 2).Value = 0.4737449735126802
$ws.Cells.Item(9, 3).Value = 0.2098342360226653
$ws.Cells.Item(9, 5).Value = 0.1367382886991191
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.5374302651227723
$ws.Cells.Item(9, 8).Value = 0.6935088678420556
$ws.Cells.Item(9, 9).Value = 0.7156976035567801
$ws.Cells.Item(9, 11).Value = 0.3385122572830994
$ws.Cells.Item(9, 12).Value = 0.2112911868116285
$ws.Cells.Item(9, 13).Value = 0.1417891887627363
$ws.Cells.Item(9, 15).Value = 2.432701736300288
$ws.Cells.Item(10, 2).Value = 0.5444910930362425
$ws.Cells.Item(10, 3).Value = 0.2100371347490082
$ws.Cells.Item(10, 5).Value = 0.1357737773411625
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.5269175336807592
$ws.Cells.Item(10, 8).Value = 0.6832746411502981
$ws.Cells.Item(10, 9).Value = 0.7020105592390191
$ws.Cells.Item(10, 11).Value = 0.3996794748327375
$ws.Cells.Item(10, 12).Value = 0.2168539693786897
$ws.Cells.Item(10, 13).Value = 0.1556703787710774
$ws.Cells.Item(10, 15).Value = 2.389308806046401
$ws.Cells.Item(11, 2).Value = 0.5766444617557624
$ws.Cells.Item(11, 3).Value = 0.2101641686234714
$ws.Cells.Item(11, 5).Value = 0.1353977639388599
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.5225110149544818
$ws.Cells.Item(11, 8).Value = 0.6789095392509523
$ws.Cells.Item(11, 9).Value = 0.6961610023871287
$ws.Cells.Item(11, 11).Value = 0.4273997026650136
$ws.Cells.Item(11, 12).Value = 0.2194855854569653
$ws.Cells.Item(11, 13).Value = 0.1620127189551823
$ws.Cells.Item(11, 15).Value = 2.370966714064338
$ws.Cells.Item(12, 2).Value = 0.5888151961987944
$ws.Cells.Item(12, 3).Value = 0.2102172253885541
$ws.Cells.Item(12, 5).Value = 0.1352643764731063
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.5208964599873624
$ws.Cells.Item(12, 8).Value = 0.6772983104987418
$ws.Cells.Item(12, 9).Value = 0.69400008459675
$ws.Cells.Item(12, 11).Value = 0.4378809750062942
$ws.Cells.Item(12, 12).Value = 0.2204965676772588
$ws.Cells.Item(12, 13).Value = 0.1644182457820591
$ws.Cells.Item(12, 15).Value = 2.364221997501261
$ws.Cells.Item(13, 2).Value = 0.5861942483632561
$ws.Cells.Item(13, 3).Value = 0.2102055791126674
$ws.Cells.Item(13, 5).Value = 0.1352927039012037
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.5212417762832189
$ws.Cells.Item(13, 8).Value = 0.6776434615346432
$ws.Cells.Item(13, 9).Value = 0.6944630670252518
$ws.Cells.Item(13, 11).Value = 0.4356243577894645
$ws.Cells.Item(13, 12).Value = 0.2202781934767017
$ws.Cells.Item(13, 13).Value = 0.1639000052658233
$ws.Cells.Item(13, 15).Value = 2.3656656507261
$ws.Cells.Item(14, 2).Value = 0.5776458612793363
$ws.Cells.Item(14, 3).Value = 0.2101684346367847
$ws.Cells.Item(14, 5).Value = 0.135386609813736
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.522377099987402
$ws.Cells.Item(14, 8).Value = 0.6787761461350712
$ws.Cells.Item(14, 9).Value = 0.6959821364460286
$ws.Cells.Item(14, 11).Value = 0.4282623232649314
$ws.Cells.Item(14, 12).Value = 0.2195684705552878
$ws.Cells.Item(14, 13).Value = 0.1622105474786508
$ws.Cells.Item(14, 15).Value = 2.370407792740735
$ws.Cells.Item(15, 2).Value = 0.5724090455570945
$ws.Cells.Item(15, 3).Value = 0.210146326174943
$ws.Cells.Item(15, 5).Value = 0.135445301381834
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.5230795655287679
$ws.Cells.Item(15, 8).Value = 0.6794753830476381
$ws.Cells.Item(15, 9).Value = 0.6969196666740736
$ws.Cells.Item(15, 11).Value = 0.4237507901253821
$ws.Cells.Item(15, 12).Value = 0.2191356234767028
$ws.Cells.Item(15, 13).Value = 0.1611761984029698
$ws.Cells.Item(15, 15).Value = 2.373338674388648
$ws.Cells.Item(16, 2).Value = 0.542389132992497
$ws.Cells.Item(16, 3).Value = 0.2100295280509172
$ws.Cells.Item(16, 5).Value = 0.1357996114316737
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.527213076919125
$ws.Cells.Item(16, 8).Value = 0.6835657535993036
$ws.Cells.Item(16, 9).Value = 0.7024004232297685
$ws.Cells.Item(16, 11).Value = 0.3978657215644716
$ws.Cells.Item(16, 12).Value = 0.2166840134640893
$ws.Cells.Item(16, 13).Value = 0.1552564355180905
$ws.Cells.Item(16, 15).Value = 2.390535633295528
$ws.Cells.Item(17, 2).Value = 0.5239647749233711
$ws.Cells.Item(17, 3).Value = 0.2099667417986382
$ws.Cells.Item(17, 5).Value = 0.1360330244323684
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.5298451416973862
$ws.Cells.Item(17, 8).Value = 0.6861494438071176
$ws.Cells.Item(17, 9).Value = 0.7058591887431778
$ws.Cells.Item(17, 11).Value = 0.3819586876268772
$ws.Cells.Item(17, 12).Value = 0.2152058557434628
$ws.Cells.Item(17, 13).Value = 0.1516318339965039
$ws.Cells.Item(17, 15).Value = 2.401443409452767
$ws.Cells.Item(18, 2).Value = 0.5133648535566522
$ws.Cells.Item(18, 3).Value = 0.2099339007039163
$ws.Cells.Item(18, 5).Value = 0.1361731848642282
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.5313944018499939
$ws.Cells.Item(18, 8).Value = 0.6876628591788645
$ws.Cells.Item(18, 9).Value = 0.7078840401274924
$ws.Cells.Item(18, 11).Value = 0.372799534951298
$ws.Cells.Item(18, 12).Value = 0.2143651778078919
$ws.Cells.Item(18, 13).Value = 0.1495496778246945
$ws.Cells.Item(18, 15).Value = 2.407848815727732
$ws.Cells.Item(19, 2).Value = 0.5097754605711486
$ws.Cells.Item(19, 3).Value = 0.209923344538943
$ws.Cells.Item(19, 5).Value = 0.1362216560693774
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.5319250269818099
$ws.Cells.Item(19, 8).Value = 0.6881799731722325
$ws.Cells.Item(19, 9).Value = 0.7085757102676649
$ws.Cells.Item(19, 11).Value = 0.3696967323904516
$ws.Cells.Item(19, 12).Value = 0.2140821763475032
$ws.Cells.Item(19, 13).Value = 0.1488451503395183
$ws.Cells.Item(19, 15).Value = 2.41004016663868
$ws.Cells.Item(20, 2).Value = 0.5259263665423362
$ws.Cells.Item(20, 3).Value = 0.209973087184828
$ws.Cells.Item(20, 5).Value = 0.1360075659836237
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.5295612931985261
$ws.Cells.Item(20, 8).Value = 0.6858715756970852
$ws.Cells.Item(20, 9).Value = 0.7054873271286155
$ws.Cells.Item(20, 11).Value = 0.3836530425872695
$ws.Cells.Item(20, 12).Value = 0.2153622235871637
$ws.Cells.Item(20, 13).Value = 0.1520174091838697
$ws.Cells.Item(20, 15).Value = 2.400268643478299
$ws.Cells.Item(21, 2).Value = 0.5801568733328395
$ws.Cells.Item(21, 3).Value = 0.2101792108217353
$ws.Cells.Item(21, 5).Value = 0.1353587832744818
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.5220421591988469
$ws.Cells.Item(21, 8).Value = 0.6784423166432916
$ws.Cells.Item(21, 9).Value = 0.6955344783869037
$ws.Cells.Item(21, 11).Value = 0.430425164764614
$ws.Cells.Item(21, 12).Value = 0.2197765420193747
$ws.Cells.Item(21, 13).Value = 0.1627066798212908
$ws.Cells.Item(21, 15).Value = 2.369009454103619
$ws.Cells.Item(22, 2).Value = 0.6155698383437027
$ws.Cells.Item(22, 3).Value = 0.210342767677254
$ws.Cells.Item(22, 5).Value = 0.1349872185254419
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.5174433235561224
$ws.Cells.Item(22, 8).Value = 0.6738301475103867
$ws.Cells.Item(22, 9).Value = 0.6893455361009018
$ws.Cells.Item(22, 11).Value = 0.4609011839835659
$ws.Cells.Item(22, 12).Value = 0.2227457172427876
$ws.Cells.Item(22, 13).Value = 0.1697149267556384
$ws.Cells.Item(22, 15).Value = 2.34975158652469
$ws.Cells.Item(23, 2).Value = 0.5966722705894654
$ws.Cells.Item(23, 3).Value = 0.2102528492642008
$ws.Cells.Item(23, 5).Value = 0.1351807376458574
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.5198689335455597
$ws.Cells.Item(23, 8).Value = 0.6762694995130403
$ws.Cells.Item(23, 9).Value = 0.6926197914030467
$ws.Cells.Item(23, 11).Value = 0.4446442224812017
$ws.Cells.Item(23, 12).Value = 0.2211533409735722
$ws.Cells.Item(23, 13).Value = 0.1659725176528255
$ws.Cells.Item(23, 15).Value = 2.359922629446146
$ws.Cells.Item(24, 2).Value = 0.5250395543093589
$ws.Cells.Item(24, 3).Value = 0.2099702082916153
$ws.Cells.Item(24, 5).Value = 0.1360190571566591
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.5296895088141369
$ws.Cells.Item(24, 8).Value = 0.6859971125952882
$ws.Cells.Item(24, 9).Value = 0.7056553325041115
$ws.Cells.Item(24, 11).Value = 0.3828870682265517
$ws.Cells.Item(24, 12).Value = 0.2152915012253516
$ws.Cells.Item(24, 13).Value = 0.1518430854101922
$ws.Cells.Item(24, 15).Value = 2.400799336598197
$ws.Cells.Item(25, 2).Value = 0.4476712205255637
$ws.Cells.Item(25, 3).Value = 0.2097905987163031
$ws.Cells.Item(25, 5).Value = 0.1371492634739511
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.5416344634913415
$ws.Cells.Item(25, 8).Value = 0.6975355310398186
$ws.Cells.Item(25, 9).Value = 0.7210720174101919
$ws.Cells.Item(25, 11).Value = 0.3158979842644101
$ws.Cells.Item(25, 12).Value = 0.2093328259504048
$ws.Cells.Item(25, 13).Value = 0.1367029384906076
$ws.Cells.Item(25, 15).Value = 2.449920386329026
